$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.246.03"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +5.26%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.915.84"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +5.76%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "253.86"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5142"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.58%  "

$ws.Range("E8").Value = "  +6.72%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.2974"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.54%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06821"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.92%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.914.10"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.17%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "17.46"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.34%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07354"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.6905"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.60%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "87.84"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.28%  "

$ws.Range("E16").Value = "  +4.43%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.250.04"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.42%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008007"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +8.46%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.03"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.06%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.162.48"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.19%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9992"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.856"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.92%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.737"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.86%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.198"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.39%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "146.34"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.61%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "139.01"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +23.17%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.30"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +7.83%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.012"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.98%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.384"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.43%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.276"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.47%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.08850"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.84%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.027"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.08%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.05136"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.33%  "

$ws.Range("E35").Value = "  +6.50%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7185"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +6.32%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.686"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.824"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.33%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.304"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.99%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.9774"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.20%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.01708"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.49%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.117"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.25%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "106.19"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.13%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.4318"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.05%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "7.697"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.77%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1278"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.46%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.05739"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.41%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "33.43"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.39%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.512"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.82%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3837"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.02%  "
